# Auto-generated edit script applying the Masamune_Profits.xlsx diff
# to the corresponding sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 2062.4285
$ws.Range("J112").Value = 2777.5715
$ws.Range("L112").Value = 8332.7145
$ws.Range("N112").Value = -10548.7145

# Row 137
$ws.Range("H137").Value = 4239.9165
$ws.Range("I137").Value = 1074.1666
$ws.Range("J137").Value = 7405.6665
$ws.Range("K137").Value = 3222.4998
$ws.Range("L137").Value = 22216.9995
$ws.Range("M137").Value = -672.4998000000001
$ws.Range("N137").Value = -27316.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9364.700999999999
$ws.Range("I32").Value = 8902.982
$ws.Range("K32").Value = 8902.982
$ws.Range("M32").Value = -8615.982

# Row 61
$ws.Range("H61").Value = 1632.5834
$ws.Range("I61").Value = 1222.7646
$ws.Range("J61").Value = 2627.8572
$ws.Range("K61").Value = 1222.7646
$ws.Range("L61").Value = 2627.8572
$ws.Range("M61").Value = -1010.7646
$ws.Range("N61").Value = -3051.8572

# Row 74
$ws.Range("H74").Value = 1950.0857
$ws.Range("I74").Value = 1752.875
$ws.Range("J74").Value = 2380.3635
$ws.Range("K74").Value = 1752.875
$ws.Range("L74").Value = 2380.3635
$ws.Range("M74").Value = -878.875
$ws.Range("N74").Value = -4128.363499999999

# Row 77
$ws.Range("H77").Value = 1950.0857
$ws.Range("I77").Value = 1752.875
$ws.Range("J77").Value = 2380.3635
$ws.Range("K77").Value = 8764.375
$ws.Range("L77").Value = 11901.8175
$ws.Range("M77").Value = -4396.375
$ws.Range("N77").Value = -20637.8175

# Row 88
$ws.Range("H88").Value = 17308956
$ws.Range("I88").Value = 50002500
$ws.Range("J88").Value = 4231537
$ws.Range("K88").Value = 50002500
$ws.Range("L88").Value = 4231537
$ws.Range("M88").Value = -50002094
$ws.Range("N88").Value = -4232349

# Row 91
$ws.Range("H91").Value = 17308956
$ws.Range("I91").Value = 50002500
$ws.Range("J91").Value = 4231537
$ws.Range("K91").Value = 50002500
$ws.Range("L91").Value = 4231537
$ws.Range("M91").Value = -50001096
$ws.Range("N91").Value = -4234345

# Row 97
$ws.Range("H97").Value = 1124.8
$ws.Range("I97").Value = 1109.1666
$ws.Range("K97").Value = 1109.1666
$ws.Range("M97").Value = -613.1666

# Row 132
$ws.Range("H132").Value = 2456.8333
$ws.Range("I132").Value = 1369.2222
$ws.Range("J132").Value = 4088.25
$ws.Range("K132").Value = 4107.6666
$ws.Range("L132").Value = 12264.75
$ws.Range("M132").Value = -1577.6666
$ws.Range("N132").Value = -17324.75

# Row 136
$ws.Range("H136").Value = 1632.5834
$ws.Range("I136").Value = 1222.7646
$ws.Range("J136").Value = 2627.8572
$ws.Range("K136").Value = 3668.2938
$ws.Range("L136").Value = 7883.571599999999
$ws.Range("M136").Value = -1118.2938
$ws.Range("N136").Value = -12983.5716

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 589933.1
$ws.Range("I86").Value = 1895.1111
$ws.Range("J86").Value = 1251475.9
$ws.Range("K86").Value = 1895.1111
$ws.Range("L86").Value = 1251475.9
$ws.Range("M86").Value = -772.1111000000001
$ws.Range("N86").Value = -1253721.9

# Row 89
$ws.Range("H89").Value = 589933.1
$ws.Range("I89").Value = 1895.1111
$ws.Range("J89").Value = 1251475.9
$ws.Range("K89").Value = 9475.5555
$ws.Range("L89").Value = 6257379.5
$ws.Range("M89").Value = -3859.5555
$ws.Range("N89").Value = -6268611.5

# Row 134
$ws.Range("H134").Value = 3000.5715
$ws.Range("I134").Value = 2331.8
$ws.Range("J134").Value = 3311.628
$ws.Range("K134").Value = 6995.400000000001
$ws.Range("L134").Value = 9934.884
$ws.Range("M134").Value = -4460.400000000001
$ws.Range("N134").Value = -15004.884

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4160.297
$ws.Range("I31").Value = 1680.8077
$ws.Range("J31").Value = 5152.0923
$ws.Range("K31").Value = 1680.8077
$ws.Range("L31").Value = 5152.0923
$ws.Range("M31").Value = -1385.8077
$ws.Range("N31").Value = -5742.0923

# Row 34
$ws.Range("H34").Value = 4160.297
$ws.Range("I34").Value = 1680.8077
$ws.Range("J34").Value = 5152.0923
$ws.Range("K34").Value = 1680.8077
$ws.Range("L34").Value = 5152.0923
$ws.Range("M34").Value = -1478.8077
$ws.Range("N34").Value = -5556.0923

# Row 58
$ws.Range("H58").Value = 1816.5667
$ws.Range("I58").Value = 1637.6
$ws.Range("J58").Value = 2174.5
$ws.Range("K58").Value = 1637.6
$ws.Range("L58").Value = 2174.5
$ws.Range("M58").Value = -1434.6
$ws.Range("N58").Value = -2580.5

# Row 132
$ws.Range("H132").Value = 46208.438
$ws.Range("I132").Value = 1265.7727
$ws.Range("J132").Value = 145082.3
$ws.Range("K132").Value = 3797.3181
$ws.Range("L132").Value = 435246.9
$ws.Range("M132").Value = -1267.3181
$ws.Range("N132").Value = -440306.9

# Row 134
$ws.Range("H134").Value = 390439.4
$ws.Range("I134").Value = 1153.8572
$ws.Range("J134").Value = 1752938.9
$ws.Range("K134").Value = 3461.5716
$ws.Range("L134").Value = 5258816.699999999
$ws.Range("M134").Value = -926.5715999999998
$ws.Range("N134").Value = -5263886.699999999

# Row 136
$ws.Range("H136").Value = 1816.5667
$ws.Range("I136").Value = 1637.6
$ws.Range("J136").Value = 2174.5
$ws.Range("K136").Value = 4912.799999999999
$ws.Range("L136").Value = 6523.5
$ws.Range("M136").Value = -2362.799999999999
$ws.Range("N136").Value = -11623.5

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 5148.4443
$ws.Range("I113").Value = 6122.278
$ws.Range("J113").Value = 3200.7778
$ws.Range("K113").Value = 18366.834
$ws.Range("L113").Value = 9602.3334
$ws.Range("M113").Value = -16196.834
$ws.Range("N113").Value = -13942.3334

# Row 131
$ws.Range("H131").Value = 3318.6875
$ws.Range("I131").Value = 8100.5386
$ws.Range("J131").Value = 1542.5714
$ws.Range("K131").Value = 24301.6158
$ws.Range("L131").Value = 4627.7142
$ws.Range("M131").Value = -19261.6158
$ws.Range("N131").Value = -14707.7142

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2484.1667
$ws.Range("I126").Value = 2720.5557
$ws.Range("J126").Value = 1775
$ws.Range("K126").Value = 8161.6671
$ws.Range("L126").Value = 5325
$ws.Range("M126").Value = -5691.6671
$ws.Range("N126").Value = -10265

# Row 132
$ws.Range("H132").Value = 2331.8206
$ws.Range("I132").Value = 1450.4762
$ws.Range("J132").Value = 3360.0557
$ws.Range("K132").Value = 4351.4286
$ws.Range("L132").Value = 10080.1671
$ws.Range("M132").Value = -1821.4286
$ws.Range("N132").Value = -15140.1671

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2204.182
$ws.Range("I132").Value = 1660.849
$ws.Range("J132").Value = 3404.0417
$ws.Range("K132").Value = 4982.547
$ws.Range("L132").Value = 10212.1251
$ws.Range("M132").Value = -2452.547
$ws.Range("N132").Value = -15272.1251

# Row 136
$ws.Range("H136").Value = 1980.3334
$ws.Range("I136").Value = 1619.5264
$ws.Range("J136").Value = 2837.25
$ws.Range("K136").Value = 4858.5792
$ws.Range("L136").Value = 8511.75
$ws.Range("M136").Value = -2308.5792
$ws.Range("N136").Value = -13611.75

$ws = $wb.Worksheets.Item("WVR")
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 132
$ws.Range("H132").Value = 1659.3954
$ws.Range("I132").Value = 1417.3103
$ws.Range("J132").Value = 2160.8572
$ws.Range("K132").Value = 4251.9309
$ws.Range("L132").Value = 6482.571599999999
$ws.Range("M132").Value = -1721.9309
$ws.Range("N132").Value = -11542.5716

# Row 136
$ws.Range("H136").Value = 334248.16
$ws.Range("I136").Value = 417335.84
$ws.Range("J136").Value = 1897.5
$ws.Range("K136").Value = 1252007.52
$ws.Range("L136").Value = 5692.5
$ws.Range("M136").Value = -1249457.52
$ws.Range("N136").Value = -10792.5
